$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 339 (pushes the existing rows 339:385 down to 340:386,
# extending the used range to A1:R386).
$ws.Rows.Item(339).Insert()

# Populate the newly inserted row 339 with the new price-observation record
# (Camote, "1a (cosecha)", 2023-02-27 / serial 44984).
$ws.Range("A339").Value = 5
$ws.Range("B339").Value = "Macroferia Regional de Talca"
$ws.Range("C339").Value = "Maule"
$ws.Range("D339").Value = 44984
$ws.Range("E339").Value = 7
$ws.Range("F339").Value = 100112045
$ws.Range("G339").Value = "Zapallo"
$ws.Range("H339").Value = "Camote"
$ws.Range("I339").Value = "1a (cosecha)"
$ws.Range("J339").Value = 900
$ws.Range("K339").Value = 300
$ws.Range("L339").Value = 300
$ws.Range("M339").Value = 300
$ws.Range("N339").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O339").Value = "Región del Maule"
$ws.Range("P339").Value = 300
$ws.Range("Q339").Value = 1
$ws.Range("R339").Value = "Hortaliza"
